$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 8523
$ws.Range("I33").Value = 8523
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 8523
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -8294
$ws.Range("N33").ClearContents()

$ws.Range("H40").Value = 1734.2222
$ws.Range("I40").Value = 1748
$ws.Range("J40").Value = 1706.6666
$ws.Range("K40").Value = 1748
$ws.Range("L40").Value = 1706.6666
$ws.Range("M40").Value = -1573
$ws.Range("N40").Value = -2056.6666

$ws.Range("I41").Value = 513.1667
$ws.Range("J41").Value = 199.66667
$ws.Range("K41").Value = 513.1667
$ws.Range("L41").Value = 199.66667
$ws.Range("M41").Value = -73.16669999999999
$ws.Range("N41").Value = -1079.66667

$ws.Range("H53").Value = 232.4
$ws.Range("I53").Value = 248
$ws.Range("J53").Value = 213.33333
$ws.Range("K53").Value = 248
$ws.Range("L53").Value = 213.33333
$ws.Range("M53").Value = 389
$ws.Range("N53").Value = -1487.33333

$ws.Range("H62").Value = 2366.7715
$ws.Range("I62").Value = 2010.6522
$ws.Range("J62").Value = 3049.3333
$ws.Range("K62").Value = 2010.6522
$ws.Range("L62").Value = 3049.3333
$ws.Range("M62").Value = -1386.6522
$ws.Range("N62").Value = -4297.3333

$ws.Range("H65").Value = 2366.7715
$ws.Range("I65").Value = 2010.6522
$ws.Range("J65").Value = 3049.3333
$ws.Range("K65").Value = 10053.261
$ws.Range("L65").Value = 15246.6665
$ws.Range("M65").Value = -6933.261
$ws.Range("N65").Value = -21486.6665

$ws.Range("H94").Value = 3333.3333
$ws.Range("J94").Value = 3000
$ws.Range("L94").Value = 3000
$ws.Range("N94").Value = -3902

$ws.Range("H98").Value = 5000
$ws.Range("I98").Value = 5000
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 5000
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -3502
$ws.Range("N98").ClearContents()

$ws.Range("H122").Value = 5000
$ws.Range("I122").Value = 5000
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 15000
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -12550
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6126.702
$ws.Range("I32").Value = 3371.0466
$ws.Range("J32").Value = 35750
$ws.Range("K32").Value = 3371.0466
$ws.Range("L32").Value = 35750
$ws.Range("M32").Value = -3084.0466
$ws.Range("N32").Value = -36324

$ws.Range("H45").Value = 2204.8
$ws.Range("I45").Value = 1512
$ws.Range("J45").Value = 2666.6667
$ws.Range("K45").Value = 1512
$ws.Range("L45").Value = 2666.6667
$ws.Range("M45").Value = -1135
$ws.Range("N45").Value = -3420.6667

$ws.Range("H74").Value = 1546.4
$ws.Range("I74").Value = 1465.7693
$ws.Range("J74").Value = 1779.3334
$ws.Range("K74").Value = 1465.7693
$ws.Range("L74").Value = 1779.3334
$ws.Range("M74").Value = -591.7692999999999
$ws.Range("N74").Value = -3527.3334

$ws.Range("H77").Value = 1546.4
$ws.Range("I77").Value = 1465.7693
$ws.Range("J77").Value = 1779.3334
$ws.Range("K77").Value = 7328.8465
$ws.Range("L77").Value = 8896.666999999999
$ws.Range("M77").Value = -2960.8465
$ws.Range("N77").Value = -17632.667

$ws.Range("H97").Value = 1433.258
$ws.Range("I97").Value = 1294.2609
$ws.Range("K97").Value = 1294.2609
$ws.Range("M97").Value = -798.2609

$ws.Range("H107").Value = 12000
$ws.Range("J107").Value = 12000
$ws.Range("L107").Value = 12000
$ws.Range("N107").Value = -19680

$ws.Range("H110").Value = 39622
$ws.Range("I110").Value = 63462.125
$ws.Range("J110").Value = 1477.8
$ws.Range("K110").Value = 63462.125
$ws.Range("L110").Value = 1477.8
$ws.Range("M110").Value = -61417.125
$ws.Range("N110").Value = -5567.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 11249
$ws.Range("I22").Value = 14437.286
$ws.Range("J22").Value = 90
$ws.Range("K22").Value = 14437.286
$ws.Range("L22").Value = 90
$ws.Range("M22").Value = -14264.286
$ws.Range("N22").Value = -436

$ws.Range("H25").Value = 4107.5
$ws.Range("I25").Value = 1707
$ws.Range("J25").Value = 6508
$ws.Range("K25").Value = 1707
$ws.Range("L25").Value = 6508
$ws.Range("M25").Value = -1472
$ws.Range("N25").Value = -6978

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H95").Value = 12766.667
$ws.Range("J95").Value = 12766.667
$ws.Range("L95").Value = 12766.667
$ws.Range("N95").Value = -18258.667

$ws.Range("H100").Value = 35189.75
$ws.Range("J100").Value = 35189.75
$ws.Range("L100").Value = 35189.75
$ws.Range("N100").Value = -37353.75

$ws.Range("H116").Value = 35742
$ws.Range("J116").Value = 35742
$ws.Range("L116").Value = 35742
$ws.Range("N116").Value = -44920

$ws.Range("H119").Value = 38840.332
$ws.Range("J119").Value = 38840.332
$ws.Range("L119").Value = 38840.332
$ws.Range("N119").Value = -48516.332

$ws.Range("H132").Value = 2309.6287
$ws.Range("I132").Value = 1626.25
$ws.Range("J132").Value = 5043.143
$ws.Range("K132").Value = 4878.75
$ws.Range("L132").Value = 15129.429
$ws.Range("M132").Value = -2348.75
$ws.Range("N132").Value = -20189.429

$ws.Range("H140").Value = 31760
$ws.Range("J140").Value = 31760
$ws.Range("L140").Value = 31760
$ws.Range("N140").Value = -42120

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 840
$ws.Range("I132").Value = 811.1111
$ws.Range("J132").Value = 1100
$ws.Range("K132").Value = 7299.9999
$ws.Range("L132").Value = 9900
$ws.Range("M132").Value = -4769.9999
$ws.Range("N132").Value = -14960

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2498.75
$ws.Range("I80").Value = 1997.5
$ws.Range("K80").Value = 1997.5
$ws.Range("M80").Value = -999.5

$ws.Range("H83").Value = 2498.75
$ws.Range("I83").Value = 1997.5
$ws.Range("K83").Value = 9987.5
$ws.Range("M83").Value = -4995.5

$ws.Range("H97").Value = 5328.5
$ws.Range("I97").Value = 5240
$ws.Range("J97").Value = 5505.5
$ws.Range("K97").Value = 5240
$ws.Range("L97").Value = 5505.5
$ws.Range("M97").Value = -4744
$ws.Range("N97").Value = -6497.5

$ws.Range("H122").Value = 2251.389
$ws.Range("I122").Value = 1936.76
$ws.Range("J122").Value = 2966.4546
$ws.Range("K122").Value = 5810.28
$ws.Range("L122").Value = 8899.363799999999
$ws.Range("M122").Value = -3360.28
$ws.Range("N122").Value = -13799.3638

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2277.2222
$ws.Range("I61").Value = 2081
$ws.Range("J61").Value = 2522.5
$ws.Range("K61").Value = 2081
$ws.Range("L61").Value = 2522.5
$ws.Range("M61").Value = -1879
$ws.Range("N61").Value = -2926.5

$ws.Range("H82").Value = 2900.6667
$ws.Range("I82").Value = 3101
$ws.Range("K82").Value = 3101
$ws.Range("M82").Value = -2740

$ws.Range("H85").Value = 2900.6667
$ws.Range("I85").Value = 3101
$ws.Range("K85").Value = 3101
$ws.Range("M85").Value = -1853

$ws.Range("H93").Value = 13488.25
$ws.Range("I93").Value = 17617.666
$ws.Range("K93").Value = 17617.666
$ws.Range("M93").Value = -16369.666

$ws.Range("H113").Value = 2277.2222
$ws.Range("I113").Value = 2081
$ws.Range("J113").Value = 2522.5
$ws.Range("K113").Value = 2081
$ws.Range("L113").Value = 2522.5
$ws.Range("M113").Value = 89
$ws.Range("N113").Value = -6862.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1355.25
$ws.Range("I96").Value = 1243
$ws.Range("J96").Value = 1435.4286
$ws.Range("K96").Value = 1243
$ws.Range("L96").Value = 1435.4286
$ws.Range("M96").Value = 130
$ws.Range("N96").Value = -4181.4286

$ws.Range("H100").Value = 1000
$ws.Range("I100").Value = 1000
$ws.Range("K100").Value = 2000
$ws.Range("M100").Value = -1459

$ws.Range("H126").Value = 37505.07
$ws.Range("I126").Value = 45108.39
$ws.Range("K126").Value = 135325.17
$ws.Range("M126").Value = -132855.17

Write-Output "Applied scheduled runner updates to Sheets"
